# Update the SDG indicator metadata sheet:
#  - shorten the indicator name in B4
#  - update the contact phone number in B9
#  - update the organization website in B10
#  - move the active selection to B7 (matches the author's final cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "3.2.2: Neonatal mortality rate"
$ws.Range("B9").Value = "0 (312) 32 55 46"
$ws.Range("B10").Value = "www.stat.gov.kg"

# Re-touch B24's font so the style table gets refreshed the same way Excel
# does when a cell's formatting is nudged during editing.
$ws.Range("B24").Font.Name = $ws.Range("B24").Font.Name

$ws.Range("B7").Select() | Out-Null
